$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-10-31 Friday"; New = "2025-11-01 Saturday" },
    @{ Old = "788×3="; New = "335×9=" },
    @{ Old = "308×5="; New = "131×5=" },
    @{ Old = "634×5="; New = "895×2=" },
    @{ Old = "170×3="; New = "420×2=" },
    @{ Old = "137×6="; New = "205×8=" },
    @{ Old = "831×2="; New = "686×3=" },
    @{ Old = "122×3="; New = "704×9=" },
    @{ Old = "304×3="; New = "403×2=" },
    @{ Old = "581×4="; New = "917×4=" },
    @{ Old = "991×6="; New = "654×9=" },
    @{ Old = "928×2="; New = "570×4=" },
    @{ Old = "346×7="; New = "840×5=" },
    @{ Old = "402×5="; New = "189×3=" },
    @{ Old = "259×3="; New = "588×8=" },
    @{ Old = "378×7="; New = "262×6=" },
    @{ Old = "840×3="; New = "181×6=" },
    @{ Old = "205×5="; New = "578×9=" },
    @{ Old = "230×7="; New = "669×4=" },
    @{ Old = "299×7="; New = "790×3=" },
    @{ Old = "380×4="; New = "411×9=" },
    @{ Old = "203×4="; New = "323×9=" },
    @{ Old = "968×4="; New = "818×7=" },
    @{ Old = "798×9="; New = "792×6=" },
    @{ Old = "901×3="; New = "820×3=" },
    @{ Old = "958×8="; New = "686×9=" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
